$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F303").Value = 9597
$ws.Range("G303").Value = 613
$ws.Range("F304").Value = 6052
$ws.Range("G304").Value = 521
$ws.Range("F305").Value = 3378
$ws.Range("G305").Value = 263
$ws.Range("F306").Value = 74926
$ws.Range("G306").Value = 7614
$ws.Range("F307").Value = 75895
$ws.Range("G307").Value = 6399
$ws.Range("F308").Value = 15471
$ws.Range("G308").Value = 1050
$ws.Range("F309").Value = 77915
$ws.Range("G309").Value = 5531
$ws.Range("F310").Value = 79229
$ws.Range("G310").Value = 4066
$ws.Range("F311").Value = 61506
$ws.Range("G311").Value = 1928
$ws.Range("F312").Value = 28135
$ws.Range("G312").Value = 926
$ws.Range("F313").Value = 75550
$ws.Range("G313").Value = 3456
$ws.Range("F314").Value = 64368
$ws.Range("G314").Value = 3149
$ws.Range("F315").Value = 56369
$ws.Range("G315").Value = 2628
$ws.Range("F316").Value = 50750
$ws.Range("G316").Value = 2299
$ws.Range("F317").Value = 63739
$ws.Range("G317").Value = 2173
$ws.Range("F318").Value = 48964
$ws.Range("G318").Value = 1135
$ws.Range("F319").Value = 41322
$ws.Range("G319").Value = 1628
$ws.Range("F320").Value = 71593
$ws.Range("G320").Value = 3305
$ws.Range("F321").Value = 89350
$ws.Range("G321").Value = 2655
$ws.Range("F322").Value = 109660
$ws.Range("G322").Value = 2338
$ws.Range("F323").Value = 216839
$ws.Range("G323").Value = 3111
$ws.Range("F324").Value = 241037
$ws.Range("G324").Value = 2809
$ws.Range("F325").Value = 766046
$ws.Range("G325").Value = 6456
$ws.Range("F326").Value = 419556
$ws.Range("G326").Value = 3835
$ws.Range("F327").Value = 225194
$ws.Range("G327").Value = 2725
$ws.Range("F328").Value = 180851
$ws.Range("G328").Value = 2671
$ws.Range("F329").Value = 82999
$ws.Range("G329").Value = 1759
$ws.Range("F330").Value = 72552
$ws.Range("G330").Value = 2084
$ws.Range("F331").Value = 155002
$ws.Range("G331").Value = 2709
$ws.Range("F332").Value = 457208
$ws.Range("G332").Value = 4535
$ws.Range("F333").Value = 271711
$ws.Range("G333").Value = 2945
